$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark Price (column D) cells as Text so numeric-looking strings
# like "63.157.65" or "1.00" are preserved verbatim instead of being
# reinterpreted as numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D12", "D13", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Step 2: write the refreshed coin data.
$ws.Range("D2").Value = "63.157.65"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.021.67"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("D5").Value = "558.71"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "154.96"
$ws.Range("E6").Value = "  -3.86%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").Value = "3.027.94"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").Value = "0.113"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("D12").Value = "0.367"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").Value = "3.544.82"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "63.210.27"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "24.09"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "3.032.70"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "397.54"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "5.10"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "11.98"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").Value = "6.65"
$ws.Range("E22").Value = "  -5.55%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "65.37"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "0.467"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").Value = "0.0₃0983"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").Value = "8.69"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").Value = "20.43"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").Value = "159.20"
$ws.Range("E33").Value = "  +4.22%  "
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("D37").Value = "1.31"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "2.533.93"
$ws.Range("E38").Value = "  -6.24%  "
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").Value = "22.81"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "3.94"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").Value = "37.48"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "0.670"
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("D44").Value = "0.0601"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "0.0249"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "5.06"
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("D48").Value = "20.15"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("D49").Value = "267.88"
$ws.Range("E49").Value = "  -3.83%  "
$ws.Range("D50").Value = "0.0949"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").Value = "10.47"
$ws.Range("E51").Value = "  -0.04%  "

# Step 3: restore the default cell style on the Price cells (the text
# formatting above was only needed to protect the literal text while
# writing; visually/structurally these cells should stay "Normal").
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
